$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Select the cell that is being edited (mirrors the selection recorded in the file)
$ws.Activate()
$ws.Range("E8").Select()

# Update the greeting text for rule R10 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"
